$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 1.726
$ws.Range("G2").Value = -0.04166666666666666
$ws.Range("H2").Value = -0.2274305555555556
$ws.Range("I2").Value = -0.1499473603888499
$ws.Range("J2").Value = -0.1499473603888499
$ws.Range("K2").Value = -14
$ws.Range("L2").Value = -2.430555555555556
$ws.Range("U2").Value = 1.14
$ws.Range("V2").Value = 0.02024866785079929
$ws.Range("W2").Value = -0.6422018348623852
$ws.Range("X2").Value = 0.09330194337382679
$ws.Range("Y2").Value = -0.735503778236212
$ws.Range("Z2").Value = 0.2704671816073499
$ws.Range("AA2").Value = -0.0405558399538338
$ws.Range("AB2").Value = 0.09304293101393067
$ws.Range("AC2").Value = -0.1335987709677645
$ws.Range("AD2").Value = 0.03
$ws.Range("AE2").Value = 0.3084839791988757
$ws.Range("AF2").Value = 0.3384839791988756
$ws.Range("AG2").Value = -0.8015160208011243
$ws.Range("AH2").Value = 0.005976218913683984
$ws.Range("AI2").Value = 0.04138712994483752
$ws.Range("AJ2").Value = -0.0144421246011249
$ws.Range("AK2").Value = -0.1138762300475327
$ws.Range("AM2").Value = -0.002
$ws.Range("AN2").Value = -0.1260504201680672
$ws.Range("AP2").Value = 3.367714373113967
$ws.Range("AQ2").Value = 452.5

# Row 3 updates
$ws.Range("D3").Value = 1.726
$ws.Range("G3").Value = -0.04166666666666666
$ws.Range("H3").Value = -0.2274305555555556
$ws.Range("I3").Value = -0.1499473603888499
$ws.Range("J3").Value = -0.1499473603888499
$ws.Range("K3").Value = -14
$ws.Range("L3").Value = -2.430555555555556
$ws.Range("U3").Value = 1.14
$ws.Range("V3").Value = 0.02024866785079929
$ws.Range("W3").Value = -0.6422018348623852
$ws.Range("X3").Value = 0.09330194337382679
$ws.Range("Y3").Value = -0.735503778236212
$ws.Range("Z3").Value = 0.2704671816073499
$ws.Range("AA3").Value = -0.0405558399538338
$ws.Range("AB3").Value = 0.09304293101393067
$ws.Range("AC3").Value = -0.1335987709677645
$ws.Range("AD3").Value = 0.03
$ws.Range("AE3").Value = 0.3084839791988757
$ws.Range("AF3").Value = 0.3384839791988756
$ws.Range("AG3").Value = -0.8015160208011243
$ws.Range("AH3").Value = 0.005976218913683984
$ws.Range("AI3").Value = 0.04138712994483752
$ws.Range("AJ3").Value = -0.0144421246011249
$ws.Range("AK3").Value = -0.1138762300475327
$ws.Range("AM3").Value = -0.002
$ws.Range("AN3").Value = -0.1260504201680672
$ws.Range("AP3").Value = 3.367714373113967
$ws.Range("AQ3").Value = 452.5

